$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the "Note" column (K)
$ws.Cells.Item(1, 11).Value2 = "Note"

# Row 56 (was "Willamette River (Marquam Brg to St. Johns Brg)") renamed
$ws.Cells.Item(56, 1).Value2 = "Willamette River (Marquam Brg to Multnomah Channel)"
$ws.Cells.Item(56, 3).Value2 = "Willamette River (Marquam Brg to Multnomah Channel)"
$ws.Cells.Item(56, 4).Value2 = "NO"
$ws.Cells.Item(56, 10).Value2 = "Willamette River (Marquam Brg to Multnomah Channel)"

# New row 57: Willamette River (Sellwood Brg to Willamette Falls)
$ws.Cells.Item(57, 1).Value2 = "Willamette River (Sellwood Brg to Willamette Falls)"
$ws.Cells.Item(57, 3).Value2 = "Willamette River (Sellwood Brg to Willamette Falls)"
$ws.Cells.Item(57, 4).Value2 = "NO"
$ws.Cells.Item(57, 10).Value2 = "Willamette River (Sellwood Brg to Willamette Falls)"

# New row 58: Lemolo Lake, requested by Hannah LaGassey (USFS) on 6/23/2022
$ws.Cells.Item(58, 1).Value2 = "*Lemolo Lake"
$ws.Cells.Item(58, 2).Value2 = "'01144938"
$ws.Cells.Item(58, 3).Value2 = "Lemolo Lake_01144938"
$ws.Cells.Item(58, 4).Value2 = "NO"
$ws.Cells.Item(58, 9).Value2 = "*Lemolo Lake_01144938"
$ws.Cells.Item(58, 10).Value2 = "Lemolo Lake_01144938"
$ws.Cells.Item(58, 11).Value2 = "Requested by Hannah LaGassey (USFS) on 6/23/2022."

# Match the author's final selection/active cell
$ws.Activate() | Out-Null
$ws.Range("J64").Select() | Out-Null
